$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "22÷6=3, 4" -> "86÷9=9, 5"
$t.Cell(1, 1).Range.Text = "86÷9=9, 5"
# Row 1, Col 2: "29÷6=4, 5" -> "88÷2=44, 0"
$t.Cell(1, 2).Range.Text = "88÷2=44, 0"
# Row 1, Col 3: "48÷7=6, 6" -> "27÷8=3, 3"
$t.Cell(1, 3).Range.Text = "27÷8=3, 3"
# Row 1, Col 4: "21÷2=10, 1" -> "85÷3=28, 1"
$t.Cell(1, 4).Range.Text = "85÷3=28, 1"
# Row 1, Col 5: "69÷2=34, 1" -> "53÷3=17, 2"
$t.Cell(1, 5).Range.Text = "53÷3=17, 2"
# Row 5, Col 1: "55÷9=6, 1" -> "70÷9=7, 7"
$t.Cell(5, 1).Range.Text = "70÷9=7, 7"
# Row 5, Col 2: "66÷3=22, 0" -> "43÷3=14, 1"
$t.Cell(5, 2).Range.Text = "43÷3=14, 1"
# Row 5, Col 3: "52÷3=17, 1" -> "93÷3=31, 0"
$t.Cell(5, 3).Range.Text = "93÷3=31, 0"
# Row 5, Col 4: "21÷2=10, 1" -> "53÷3=17, 2"
$t.Cell(5, 4).Range.Text = "53÷3=17, 2"
# Row 5, Col 5: "67÷3=22, 1" -> "37÷9=4, 1"
$t.Cell(5, 5).Range.Text = "37÷9=4, 1"
# Row 9, Col 1: "28÷4=7, 0" -> "76÷7=10, 6"
$t.Cell(9, 1).Range.Text = "76÷7=10, 6"
# Row 9, Col 2: "87÷6=14, 3" -> "68÷4=17, 0"
$t.Cell(9, 2).Range.Text = "68÷4=17, 0"
# Row 9, Col 3: "51÷6=8, 3" -> "53÷3=17, 2"
$t.Cell(9, 3).Range.Text = "53÷3=17, 2"
# Row 9, Col 4: "40÷3=13, 1" -> "69÷6=11, 3"
$t.Cell(9, 4).Range.Text = "69÷6=11, 3"
# Row 9, Col 5: "68÷2=34, 0" -> "66÷3=22, 0"
$t.Cell(9, 5).Range.Text = "66÷3=22, 0"
# Row 13, Col 1: "99÷7=14, 1" -> "57÷8=7, 1"
$t.Cell(13, 1).Range.Text = "57÷8=7, 1"
# Row 13, Col 2: "36÷8=4, 4" -> "13÷6=2, 1"
$t.Cell(13, 2).Range.Text = "13÷6=2, 1"
# Row 13, Col 3: "83÷5=16, 3" -> "81÷5=16, 1"
$t.Cell(13, 3).Range.Text = "81÷5=16, 1"
# Row 13, Col 4: "48÷7=6, 6" -> "72÷8=9, 0"
$t.Cell(13, 4).Range.Text = "72÷8=9, 0"
# Row 13, Col 5: "68÷3=22, 2" -> "53÷8=6, 5"
$t.Cell(13, 5).Range.Text = "53÷8=6, 5"
# Row 17, Col 1: "91÷3=30, 1" -> "49÷5=9, 4"
$t.Cell(17, 1).Range.Text = "49÷5=9, 4"
# Row 17, Col 2: "30÷2=15, 0" -> "53÷9=5, 8"
$t.Cell(17, 2).Range.Text = "53÷9=5, 8"
# Row 17, Col 3: "78÷7=11, 1" -> "41÷5=8, 1"
$t.Cell(17, 3).Range.Text = "41÷5=8, 1"
# Row 17, Col 4: "99÷2=49, 1" -> "90÷8=11, 2"
$t.Cell(17, 4).Range.Text = "90÷8=11, 2"
# Row 17, Col 5: "16÷8=2, 0" -> "55÷9=6, 1"
$t.Cell(17, 5).Range.Text = "55÷9=6, 1"
